$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# The panel was re-queried; refresh the "time_taken" timestamps on "data".
$dataSheet.Range("F2").Value = "2021-10-05 14:19:16.554289"
$dataSheet.Range("F3").Value = "2021-10-05 14:19:16.554302"
$dataSheet.Range("F4").Value = "2021-10-05 14:19:16.554307"

# Add a new "metadata" worksheet right after "data".
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# Header row.
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Reuse the header/index cell formatting already used on the "data" sheet.
$dataSheet.Range("B1:F1").Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)
$dataSheet.Range("B1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

$ws.Range("A2").Value = 0
$dataSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# Data row.
$ws.Range("B2").Value = "Autosomal recessive primary hypertrophic osteoarthropathy"
$ws.Range("C2").Value = 557

# Force "1.9" to be stored as text (not a number), matching the source data.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1.9"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "2021-03-15T22:27:06.013313Z"
$ws.Range("F2").Value = "2021-10-05 14:19:16.550322"
$ws.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/557/?format=json"

# Keep "data" as the active/selected tab (matches the source workbook).
$dataSheet.Activate()
